# Update the 2024 (row 26) values for "Energy Storage" (column C) and
# "Solar" (column E) with the latest upstream data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C26").Value = 27.68
$ws.Range("E26").Value = 2373.202
